$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos - daily totals refresh ---
$ws.Range("B4").Value = 1846873
$ws.Range("C4").Value = 9703
$ws.Range("D4").Value = 607249
$ws.Range("E4").Value = 1133117
$ws.Range("G4").Value = 312
$ws.Range("H4").Value = 106507

# --- Row 37: Egipto - daily totals refresh ---
$ws.Range("D37").Value = 6447
$ws.Range("E37").Value = 18932

# --- Rows 115-118: Costa Rica updates & overtakes Tunez, Republica de Africa
#     Central and Letonia in the ranking, pushing them down one place ---
$ws.Range("A115").Value = "Costa Rica"
$ws.Range("B115").Value = 1084
$ws.Range("C115").Value = 28
$ws.Range("D115").Value = 676
$ws.Range("E115").Value = 398
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 10

$ws.Range("A116").Value = "Tunez"
$ws.Range("B116").Value = 1084
$ws.Range("C116").Value = 7
$ws.Range("D116").Value = 964
$ws.Range("E116").Value = 72
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 48

$ws.Range("A117").Value = "Republica de Africa Central"
$ws.Range("B117").Value = 1069
$ws.Range("C117").Value = 58
$ws.Range("D117").Value = 23
$ws.Range("E117").Value = 1042
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 2
$ws.Range("H117").Value = 4

$ws.Range("A118").Value = "Letonia"
$ws.Range("B118").Value = 1066
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 745
$ws.Range("E118").Value = 297
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 24

# --- Rows 145-146: Togo updates & overtakes Taiwan in the ranking ---
$ws.Range("A145").Value = "Togo"
$ws.Range("B145").Value = 443
$ws.Range("C145").Value = 1
$ws.Range("D145").Value = 215
$ws.Range("E145").Value = 215
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 13

$ws.Range("A146").Value = "Taiwan"
$ws.Range("B146").Value = 443
$ws.Range("C146").Value = 1
$ws.Range("D146").Value = 427
$ws.Range("E146").Value = 9
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 7

# --- Update the "last updated" timestamp footer ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 21:35"
